$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 536.82355
$ws.Range("J2").Value = 584.7143
$ws.Range("L2").Value = 584.7143
$ws.Range("N2").Value = -810.7143
$ws.Range("H111").Value = 7799.4287
$ws.Range("I111").Value = 6919.2
$ws.Range("K111").Value = 20757.6
$ws.Range("M111").Value = -17690.6
$ws.Range("H125").Value = 4498.857
$ws.Range("I125").Value = 493
$ws.Range("J125").Value = 5166.5
$ws.Range("K125").Value = 4437
$ws.Range("L125").Value = 46498.5
$ws.Range("M125").Value = -1977
$ws.Range("N125").Value = -51418.5
$ws.Range("H135").Value = 8930076
$ws.Range("I135").Value = 1175.8422
$ws.Range("J135").Value = 27779976
$ws.Range("K135").Value = 10582.5798
$ws.Range("L135").Value = 250019784
$ws.Range("M135").Value = -8047.5798
$ws.Range("N135").Value = -250024854
$ws.Range("H138").Value = 3206.0652
$ws.Range("J138").Value = 4063.25
$ws.Range("L138").Value = 12189.75
$ws.Range("N138").Value = -22469.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19621282
$ws.Range("I32").Value = 30310588
$ws.Range("J32").Value = 24221.166
$ws.Range("K32").Value = 30310588
$ws.Range("L32").Value = 24221.166
$ws.Range("M32").Value = -30310301
$ws.Range("N32").Value = -24795.166
$ws.Range("H61").Value = 52637020
$ws.Range("I61").Value = 76926536
$ws.Range("K61").Value = 76926536
$ws.Range("M61").Value = -76926324
$ws.Range("H74").Value = 66743970
$ws.Range("I74").Value = 71510970
$ws.Range("J74").Value = 5999
$ws.Range("K74").Value = 71510970
$ws.Range("L74").Value = 5999
$ws.Range("M74").Value = -71510096
$ws.Range("N74").Value = -7747
$ws.Range("H77").Value = 66743970
$ws.Range("I77").Value = 71510970
$ws.Range("J77").Value = 5999
$ws.Range("K77").Value = 357554850
$ws.Range("L77").Value = 29995
$ws.Range("M77").Value = -357550482
$ws.Range("N77").Value = -38731
$ws.Range("H136").Value = 52637020
$ws.Range("I136").Value = 76926536
$ws.Range("K136").Value = 230779608
$ws.Range("M136").Value = -230777058
$ws.Range("H137").Value = 27498.5
$ws.Range("J137").Value = 27498.5
$ws.Range("L137").Value = 27498.5
$ws.Range("N137").Value = -37698.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 105263
$ws.Range("J60").Value = 105263
$ws.Range("L60").Value = 105263
$ws.Range("N60").Value = -106461
$ws.Range("H94").Value = 2632.375
$ws.Range("I94").Value = 877.26086
$ws.Range("J94").Value = 43000
$ws.Range("K94").Value = 877.26086
$ws.Range("L94").Value = 43000
$ws.Range("M94").Value = -426.26086
$ws.Range("N94").Value = -43902
$ws.Range("H106").Value = 12884.833
$ws.Range("J106").Value = 12884.833
$ws.Range("L106").Value = 12884.833
$ws.Range("N106").Value = -15408.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22731240
$ws.Range("I31").Value = 3195.5715
$ws.Range("K31").Value = 3195.5715
$ws.Range("M31").Value = -2900.5715
$ws.Range("H34").Value = 22731240
$ws.Range("I34").Value = 3195.5715
$ws.Range("K34").Value = 3195.5715
$ws.Range("M34").Value = -2993.5715
$ws.Range("H125").Value = 4000
$ws.Range("J125").Value = 4000
$ws.Range("L125").Value = 4000
$ws.Range("N125").Value = -8920
$ws.Range("H132").Value = 146417.22
$ws.Range("I132").Value = 184528.45
$ws.Range("K132").Value = 553585.3500000001
$ws.Range("M132").Value = -551055.3500000001
$ws.Range("H134").Value = 1403.7941
$ws.Range("I134").Value = 1246.2069
$ws.Range("J134").Value = 2317.8
$ws.Range("K134").Value = 3738.620699999999
$ws.Range("L134").Value = 6953.400000000001
$ws.Range("M134").Value = -1203.620699999999
$ws.Range("N134").Value = -12023.4
$ws.Range("H136").Value = 2367.2632
$ws.Range("I136").Value = 2131.8125
$ws.Range("J136").Value = 3623
$ws.Range("K136").Value = 6395.4375
$ws.Range("L136").Value = 10869
$ws.Range("M136").Value = -3845.4375
$ws.Range("N136").Value = -15969
$ws.Range("H141").Value = 327999
$ws.Range("J141").Value = 340234.25
$ws.Range("L141").Value = 340234.25
$ws.Range("N141").Value = -350594.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H80").Value = 4623.4
$ws.Range("I80").Value = 3801.6667
$ws.Range("J80").Value = 5171.222
$ws.Range("K80").Value = 3801.6667
$ws.Range("L80").Value = 5171.222
$ws.Range("M80").Value = -2803.6667
$ws.Range("N80").Value = -7167.222
$ws.Range("H83").Value = 4623.4
$ws.Range("I83").Value = 3801.6667
$ws.Range("J83").Value = 5171.222
$ws.Range("K83").Value = 19008.3335
$ws.Range("L83").Value = 25856.11
$ws.Range("M83").Value = -14016.3335
$ws.Range("N83").Value = -35840.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3390
$ws.Range("J22").Value = 4085.1428
$ws.Range("L22").Value = 4085.1428
$ws.Range("N22").Value = -4675.1428
$ws.Range("H27").Value = 3390
$ws.Range("J27").Value = 4085.1428
$ws.Range("L27").Value = 4085.1428
$ws.Range("N27").Value = -4299.1428
$ws.Range("H46").Value = 1326.08
$ws.Range("I46").Value = 899.94446
$ws.Range("J46").Value = 2421.8572
$ws.Range("K46").Value = 899.94446
$ws.Range("L46").Value = 2421.8572
$ws.Range("M46").Value = -711.94446
$ws.Range("N46").Value = -2797.8572
$ws.Range("H68").Value = 4717.5713
$ws.Range("I68").Value = 3111.7144
$ws.Range("J68").Value = 6323.4287
$ws.Range("K68").Value = 3111.7144
$ws.Range("L68").Value = 6323.4287
$ws.Range("M68").Value = -2362.7144
$ws.Range("N68").Value = -7821.4287
$ws.Range("H71").Value = 4717.5713
$ws.Range("I71").Value = 3111.7144
$ws.Range("J71").Value = 6323.4287
$ws.Range("K71").Value = 15558.572
$ws.Range("L71").Value = 31617.1435
$ws.Range("M71").Value = -11814.572
$ws.Range("N71").Value = -39105.14350000001
$ws.Range("H132").Value = 50003070
$ws.Range("I132").Value = 3153.652
$ws.Range("K132").Value = 9460.956
$ws.Range("M132").Value = -6930.956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 683.5
$ws.Range("I4").Value = 641
$ws.Range("J4").Value = 760
$ws.Range("K4").Value = 641
$ws.Range("L4").Value = 760
$ws.Range("M4").Value = -528
$ws.Range("N4").Value = -986
$ws.Range("H136").Value = 1520.4
$ws.Range("I136").Value = 964.5333000000001
$ws.Range("J136").Value = 3188
$ws.Range("K136").Value = 2893.5999
$ws.Range("L136").Value = 9564
$ws.Range("M136").Value = -343.5999000000002
$ws.Range("N136").Value = -14664
